$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.293.80"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "3.492.78"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'587.07"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").Value = "'134.21"
$ws.Range("E6").Value = "  +1.72%  "

$ws.Range("D7").Value = "3.491.58"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.485"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("D11").Value = "'7.19"
$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("E12").Value = "  -2.38%  "

$ws.Range("D13").Value = "4.092.13"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").Value = "'0.120"

$ws.Range("E15").Value = "  +0.18%  "

$ws.Range("D16").Value = "3.497.13"
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("D17").Value = "64.337.50"
$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("D18").Value = "'25.15"
$ws.Range("E18").Value = "  -9.44%  "

$ws.Range("D19").Value = "'9.99"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").Value = "'5.73"
$ws.Range("E20").Value = "  +1.27%  "

$ws.Range("D21").Value = "'13.63"
$ws.Range("E21").Value = "  -5.82%  "

$ws.Range("D22").Value = "'388.27"
$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("B23").Value = "WrappedeETH"
$ws.Range("C23").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D23").Value = "3.635.06"
$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.565"
$ws.Range("E24").Value = "  -2.38%  "

$ws.Range("D25").Value = "'74.64"
$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.54"
$ws.Range("E30").Value = "  -3.15%  "

$ws.Range("D31").Value = "'7.38"
$ws.Range("E31").Value = "  -1.60%  "

$ws.Range("D32").Value = "'8.27"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -1.34%  "

$ws.Range("D34").Value = "3.514.80"
$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +2.17%  "

$ws.Range("D37").Value = "'23.46"
$ws.Range("E37").Value = "  -2.05%  "

$ws.Range("D38").Value = "'5.24"
$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("D39").Value = "'6.84"
$ws.Range("E39").Value = "  -2.26%  "

$ws.Range("E40").Value = "  -2.42%  "

$ws.Range("D41").Value = "'161.39"
$ws.Range("E41").Value = "  -3.99%  "

$ws.Range("D42").Value = "'0.0779"
$ws.Range("E42").Value = "  -3.87%  "

$ws.Range("D43").Value = "'0.805"
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'25.46"
$ws.Range("E45").Value = "  -5.45%  "

$ws.Range("D46").Value = "'41.76"
$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("D47").Value = "'4.40"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("D48").Value = "'1.18"
$ws.Range("E48").Value = "  -1.23%  "

$ws.Range("D49").Value = "'1.66"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").Value = "2.471.31"
$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("D51").Value = "'6.73"
$ws.Range("E51").Value = "  -2.52%  "
